# Scen_COM_FR_ELE.xlsx - "Add files via upload"
#
# The Cset_CN ("Commodity Set: Commodity Name") column for the COM_FR_2023
# scenario rows (C7:C230) changes from "ELC_FIN" to "ELC_FIN,ELC_IND_FIN" -
# i.e. the electricity-final-demand commodity set now also includes the
# industrial-electricity-final-demand commodity.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COM_FR_2023")
$ws.Activate()

$dataRange = $ws.Range("C7:C230")
$dataRange.Value = "ELC_FIN,ELC_IND_FIN"

# Leave the selection/viewport the way the author had it when the file was
# last saved: scrolled down so row 196 is at the top, with the whole edited
# column selected (active cell on the first row of the selection).
$excel.ActiveWindow.ScrollRow = 196
$excel.ActiveWindow.ScrollColumn = 1
$dataRange.Select() | Out-Null
